$d = $word.ActiveDocument

$d.Content.Find.Execute("72×70=5040", $true, $false, $false, $false, $false, $true, 1, $false, "41×45=1845", 2) | Out-Null
$d.Content.Find.Execute("57×89=5073", $true, $false, $false, $false, $false, $true, 1, $false, "36×40=1440", 2) | Out-Null
$d.Content.Find.Execute("60×56=3360", $true, $false, $false, $false, $false, $true, 1, $false, "34×43=1462", 2) | Out-Null
$d.Content.Find.Execute("45×37=1665", $true, $false, $false, $false, $false, $true, 1, $false, "43×40=1720", 2) | Out-Null
$d.Content.Find.Execute("68×61=4148", $true, $false, $false, $false, $false, $true, 1, $false, "90×80=7200", 2) | Out-Null
$d.Content.Find.Execute("67×23=1541", $true, $false, $false, $false, $false, $true, 1, $false, "22×64=1408", 2) | Out-Null
$d.Content.Find.Execute("71×21=1491", $true, $false, $false, $false, $false, $true, 1, $false, "61×92=5612", 2) | Out-Null
$d.Content.Find.Execute("33×75=2475", $true, $false, $false, $false, $false, $true, 1, $false, "61×26=1586", 2) | Out-Null
$d.Content.Find.Execute("37×32=1184", $true, $false, $false, $false, $false, $true, 1, $false, "55×86=4730", 2) | Out-Null
$d.Content.Find.Execute("76×75=5700", $true, $false, $false, $false, $false, $true, 1, $false, "30×17=510", 2) | Out-Null
$d.Content.Find.Execute("97×81=7857", $true, $false, $false, $false, $false, $true, 1, $false, "95×46=4370", 2) | Out-Null
$d.Content.Find.Execute("67×60=4020", $true, $false, $false, $false, $false, $true, 1, $false, "15×38=570", 2) | Out-Null
$d.Content.Find.Execute("86×15=1290", $true, $false, $false, $false, $false, $true, 1, $false, "48×21=1008", 2) | Out-Null
$d.Content.Find.Execute("90×18=1620", $true, $false, $false, $false, $false, $true, 1, $false, "44×94=4136", 2) | Out-Null
$d.Content.Find.Execute("55×73=4015", $true, $false, $false, $false, $false, $true, 1, $false, "91×80=7280", 2) | Out-Null
$d.Content.Find.Execute("71×91=6461", $true, $false, $false, $false, $false, $true, 1, $false, "39×48=1872", 2) | Out-Null
$d.Content.Find.Execute("82×50=4100", $true, $false, $false, $false, $false, $true, 1, $false, "48×31=1488", 2) | Out-Null
$d.Content.Find.Execute("23×87=2001", $true, $false, $false, $false, $false, $true, 1, $false, "35×31=1085", 2) | Out-Null
$d.Content.Find.Execute("48×36=1728", $true, $false, $false, $false, $false, $true, 1, $false, "97×35=3395", 2) | Out-Null
$d.Content.Find.Execute("71×59=4189", $true, $false, $false, $false, $false, $true, 1, $false, "90×66=5940", 2) | Out-Null
$d.Content.Find.Execute("38×19=722", $true, $false, $false, $false, $false, $true, 1, $false, "34×97=3298", 2) | Out-Null
$d.Content.Find.Execute("53×26=1378", $true, $false, $false, $false, $false, $true, 1, $false, "49×45=2205", 2) | Out-Null
$d.Content.Find.Execute("51×49=2499", $true, $false, $false, $false, $false, $true, 1, $false, "92×81=7452", 2) | Out-Null
$d.Content.Find.Execute("93×95=8835", $true, $false, $false, $false, $false, $true, 1, $false, "77×78=6006", 2) | Out-Null
$d.Content.Find.Execute("20×51=1020", $true, $false, $false, $false, $false, $true, 1, $false, "14×38=532", 2) | Out-Null
